# ------------------------------------------------------------------
# Refresh the cryptos list (Coin / Link / Price / Volume(1h)) with the
# latest values pulled from coinranking.com, as performed by the
# scheduled "Updated cryptos list ... with GitHub Actions" workflow.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some new Price values look like plain numbers (e.g. "1.001"). Mark those
# cells as Text first so Excel keeps them as strings instead of silently
# converting them to numeric values (matching the original inline-string
# cells produced by the data-refresh script).
$ws.Range("D4:D12").NumberFormat = "@"
$ws.Range("D14:D17").NumberFormat = "@"
$ws.Range("D19:D21").NumberFormat = "@"
$ws.Range("D23:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.478.13"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.930.60"
$ws.Range("E3").Value = "  +4.64%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "243.17"
$ws.Range("E5").Value = "  +4.34%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.4675"
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("D8").Value = "44.48"
$ws.Range("E8").Value = "  +2.92%  "
$ws.Range("D9").Value = "0.2890"
$ws.Range("E9").Value = "  +5.19%  "
$ws.Range("D10").Value = "0.07043"
$ws.Range("E10").Value = "  +11.29%  "
$ws.Range("D11").Value = "107.32"
$ws.Range("E11").Value = "  +25.73%  "
$ws.Range("D12").Value = "18.60"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("D13").Value = "1.907.61"
$ws.Range("E13").Value = "  +5.94%  "
$ws.Range("D14").Value = "0.07660"
$ws.Range("E14").Value = "  +2.82%  "
$ws.Range("D15").Value = "5.208"
$ws.Range("E15").Value = "  +5.62%  "
$ws.Range("D16").Value = "0.6651"
$ws.Range("E16").Value = "  +6.79%  "
$ws.Range("D17").Value = "299.27"
$ws.Range("E17").Value = "  +12.69%  "
$ws.Range("D18").Value = "30.488.43"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000007740"
$ws.Range("E19").Value = "  +5.69%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "13.11"
$ws.Range("E20").Value = "  +3.54%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "2.159.49"
$ws.Range("E22").Value = "  +4.19%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.307"
$ws.Range("E23").Value = "  +7.81%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "6.288"
$ws.Range("E25").Value = "  +6.71%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "21.94"
$ws.Range("E26").Value = "  +18.73%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "168.13"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "9.300"
$ws.Range("E28").Value = "  +3.84%  "
$ws.Range("D29").Value = "2.068"
$ws.Range("E29").Value = "  +9.53%  "
$ws.Range("D30").Value = "0.1099"
$ws.Range("E30").Value = "  +7.39%  "
$ws.Range("D31").Value = "1.362"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D32").Value = "4.136"
$ws.Range("E32").Value = "  +2.68%  "
$ws.Range("D33").Value = "4.002"
$ws.Range("E33").Value = "  +4.34%  "
$ws.Range("D34").Value = "0.05185"
$ws.Range("E34").Value = "  +7.77%  "
$ws.Range("D35").Value = "0.7422"
$ws.Range("E35").Value = "  +7.33%  "
$ws.Range("D36").Value = "1.155"
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "1.000"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "2.737"
$ws.Range("E38").Value = "  +1.95%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.02033"
$ws.Range("E39").Value = "  +9.32%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.688"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "2.067"
$ws.Range("E41").Value = "  +4.60%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.8831"
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "109.59"
$ws.Range("E43").Value = "  +3.54%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.848"
$ws.Range("E44").Value = "  +7.29%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.4281"
$ws.Range("E45").Value = "  +5.69%  "
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").Value = "53.50"
$ws.Range("E46").Value = "  +31.81%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "68.89"
$ws.Range("E48").Value = "  +11.10%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "7.243"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.297"
$ws.Range("E50").Value = "  +8.54%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.1218"
$ws.Range("E51").Value = "  +2.22%  "
